$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Duplicate "Include from FSIII" sheet as "Include from FSIII 2" ---
# (the new sheet preserves the old concept code "C" that the source sheet had)
$src = $wb.Worksheets.Item("Include from FSIII")
$src.Copy($null, $src)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Include from FSIII 2"

# --- Update the original "Include from FSIII" sheet to the new concept code ---
$src.Range("C2").Value = "d6d48a71-b96f-4b88-86f9-b13bd3c03560"

# --- Keep the originally active tab selected ---
$meta.Activate()
